$d = $word.ActiveDocument

# --- Remove the four "IT Work" bullet-point questions plus the blank
#     paragraph that followed them (they sat right after the heading). ---
$delStart = $d.Paragraphs(2).Range.Start
$delEnd   = $d.Paragraphs(6).Range.End
$d.Range($delStart, $delEnd).Delete()

# --- Append the rest of the interview write-up after the final existing
#     paragraph ("Another type of testing ... stories and more."). ---
$newParas = @(
  "",
  "One of the other types of testing Shane was required to complete was testing to make sure that the installer worked properly on each platform they were launching on. He says its ‘part of the job that you really don’t hear about’. With a total of 5 computers to test, Shane was tasked with installing the game, then checking that the binary files had correctly assigned digital certificates and check that it integrated with games for Windows live and could launch from there, and that is could then also be uninstalled and left no traces of the game on the machine, which he described at quite a lengthy process.",
  "",
  "Shane explains that his role in QA didn’t actually have many interactions with other people. He explained he didn’t even really have much to do with the person who was directly above him, the senior QA, they would come around to collect the daily smoke tests and make sure the QA department had everything they needed such as the latest builds. Shane instead had more contact with the producer of the whole QA team. Every bug that was found entered into what Shane referred to as ‘triage’, this involved that when a bug was found it would be assigned to the producer to go through them and make sure they went to the right person or department, or get rid of them if there were duplicates.",
  "",
  "The other outside his own department that Shane describes having contact with was the IT department, although Shane describes him only needing to contact them when ‘shit was broken’. As well as brief contact with the Level designer and UI designer when there were some major bugs in relation to those elements and they needed a better understanding of the bug. In regards to contact with the IT department, Shane briefly mentions the time where his emails stopped working, but the main contact during his time at 2K with the IT department was to do with rain. Shane stresses that ‘this isn’t terribly reflective of a usual workplace’ but due to the location and age of the office he was working in the roof of the building collapsed a few times when heavy rain hit four or five times, Shane explains that it was the IT departments job to deal with that issue.",
  "",
  "For day to day work Shane described his work schedule consisting mostly of the morning smoke test, which when he was new to the job would take the first half of the day, then the rest of the day was spent on general bug testing. He described the general bug testing being split between finding new bugs and checking old ones that had been sent back, and confirming they had been fixed.",
  "",
  "When asked what the most challenging part of that role was he mentioned the challenges of testing the multiplayer levels of Bioshock 2. He explained that one of the issues was that the multiplayer portion of the game was developed by an entirely separate studio and described it as ‘basically a different game’ that he was disappointed by. “It was hard to test in general because it had to be integrated through games for Windows Live” the team had to mess around with dummy accounts, make sure they were all hooked up together, he also described testing the multiplayer portion as ‘really contrived’ not enjoying the multiplayer component at all.",
  "",
  "Finally when asked if he could share an example of what best captures his work at 2K on bioshock2, Shane opted to share his experience with a slightly different game, although he was at 2K working on Bioshock, it was when 2K’s sister company gearbox software brought in the original borderlands for the team to play that Shane explains how incredible it is to think how the team he was on got to shape the game he wasn’t even employed to work on. You can hear that story at the rest of the interview by visiting our website and going to the IT Interview page.",
  "",
  "I would finally like to give a big thank you to Shane for taking the time out of his busy schedule to sit down with me and speak about his time in QA at 2K. It was not only really insightful but interesting to learn about the different types of testing games studios do as well as some of the problems that do come up in development."
)

foreach ($t in $newParas) {
  $n = $d.Paragraphs.Count
  $d.Paragraphs($n).Range.InsertParagraphAfter()
  if ($t -ne "") {
    $n = $d.Paragraphs.Count
    $d.Paragraphs($n).Range.Text = $t
  }
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
